# Adds the "Week 10" section to the end of the progress log.
#
# Strategy: the document currently ends with a paragraph that contains only
# the (hidden) "_GoBack" bookmark. We need to insert a Heading1 + a bunch of
# ListParagraph bullet items before that paragraph, turn the bookmark
# paragraph itself into a bullet item (splitting its text around the
# bookmark), and then append more bullet items plus a final empty paragraph
# after it.
#
# Word's InsertXML (as implemented here) replaces the *entire* paragraph
# that contains the target range, so every insertion point below first gets
# its own dedicated, originally-empty placeholder paragraph (created with
# InsertParagraphBefore/After) before InsertXML is used to fill it in. That
# way neighbouring paragraphs - and the bookmark - are never clobbered.

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wordNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# The paragraph that (currently) only holds the "_GoBack" bookmark.
$bookmarkP = $d.Paragraphs.Last

# ---------------------------------------------------------------------
# 1) Insert the Week 10 heading + first 5 bullet items BEFORE the bookmark
#    paragraph.
# ---------------------------------------------------------------------
$bookmarkP.Range.InsertParagraphBefore()
$placeholderBefore = $d.Paragraphs($d.Paragraphs.Count - 1)

$bodyBefore = @'
<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Week 10 (23</w:t></w:r><w:r><w:t xml:space="preserve"> hours)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Updated all menu backgrounds to use the new cycling system</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added the third-party Victory plugin, allowing…</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>The player’s external IP address to be shown</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Volume sliders in the options menu</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>A full list of supported resolutions available</w:t></w:r></w:p>
'@

$placeholderBefore.Range.InsertXML((New-PkgXml $bodyBefore))

# The InsertXML call above re-seats every paragraph handle, so re-resolve
# the (untouched, still-empty) bookmark paragraph - it is simply the last
# paragraph in the document again.
$bookmarkP = $d.Paragraphs.Last

# ---------------------------------------------------------------------
# 2) Turn the bookmark paragraph into the "Added a smooth transition from…"
#    bullet, keeping the _GoBack bookmark sitting inside the run split.
# ---------------------------------------------------------------------
$bodyBookmark = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Added a smooth transition fro</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>m the location where you die to the spectator camera, waiting for a few seconds before moving the camera</w:t></w:r></w:p>'

$bookmarkP.Range.InsertXML((New-PkgXml $bodyBookmark))

# Re-resolve the paragraph that now holds the bookmark (same paragraph,
# freshly re-seated handle) so we can anchor the next insertions after it.
$bookmarkP = $d.Paragraphs.Last

# ---------------------------------------------------------------------
# 3) Insert the remaining bullet items + a trailing empty paragraph AFTER
#    the bookmark paragraph.
# ---------------------------------------------------------------------
$bookmarkP.Range.InsertParagraphAfter()
$placeholderAfter = $d.Paragraphs.Last

$bodyAfter = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Put in the full list of credits and animated them</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Small fixes and tweaks</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Fixed camera follow height in some situations</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added collision to the new chandelier</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Made the coffin closing sound louder and carry further across the map</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added some decorative trees around the outside of the map</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added sound effect to when you try to throw a bomb but have none</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Made the bone sounds a bit more </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>boney</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Made the detach-all-limbs debug key only work if you hold it down for a second (so players are less likely to accidentally use it)</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added another boundary wall to stop players getting stuck outside the fence</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Tweaked the eye adaptation tolerances so it’s less irritating</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added a sound effect to the title-fly-in animation in the main menu</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Added an option to host a game as “unlisted,” which will work even when the session nodes fail </w:t></w:r></w:p>
<w:p/>
'@

$placeholderAfter.Range.InsertXML((New-PkgXml $bodyAfter))

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
